$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 613: header for the newly appended basket block
$ws.Range("A613").Value = "Date"
$ws.Range("B613").Value = "Stocks"
$ws.Range("C613").Value = "INFY.NS"
$ws.Range("D613").Value = "TCS.NS"
$ws.Range("E613").Value = "LTIM.NS"
$ws.Range("F613").Value = "WIPRO.NS"
$ws.Range("G613").Value = "TATATECH.NS"
$ws.Range("H613").Value = "Basket Value"
$ws.Range("J613").Value = "NAV"

# Bold the header/label cells (new bold font, mirrors the other header rows)
$ws.Range("A613").Font.Bold = $true
$ws.Range("B613").Font.Bold = $true
$ws.Range("H613").Font.Bold = $true
$ws.Range("J613").Font.Bold = $true

# Row 614: Quantities
$ws.Range("B614").Value = "Quantities"
$ws.Range("B614").Font.Bold = $true
$ws.Range("C614").Value = 5.0
$ws.Range("D614").Value = 2.0
$ws.Range("E614").Value = 1.0
$ws.Range("F614").Value = 12.0
$ws.Range("G614").Value = 5.0

# Rows 615-636: daily basket data for the new basket
# Force column A to Text first so these date-like strings are stored as literal
# text (matching the source data) instead of being auto-converted to date serials
$ws.Range("A615:A636").NumberFormat = "@"

$rngData = New-Object 'object[,]' 22,10
$rngData[0,0] = "2024-08-26"
$rngData[0,2] = 1876.150024414062
$rngData[0,3] = 4502.4501953125
$rngData[0,4] = 5739.9501953125
$rngData[0,5] = 520.0
$rngData[0,6] = 1048.650024414062
$rngData[0,7] = 35608.85083007812
$rngData[0,8] = 0.0
$rngData[0,9] = 179.1548046208872
$rngData[1,0] = "2024-08-28"
$rngData[1,2] = 1900.099975585938
$rngData[1,3] = 4497.14990234375
$rngData[1,4] = 5751.5498046875
$rngData[1,5] = 517.1500244140625
$rngData[1,6] = 1078.800048828125
$rngData[1,7] = 35846.15002441406
$rngData[1,8] = 0.006664050897578962
$rngData[1,9] = 180.3487013574266
$rngData[2,0] = "2024-08-29"
$rngData[2,2] = 1939.099975585938
$rngData[2,3] = 4506.0498046875
$rngData[2,4] = 6127.5498046875
$rngData[2,5] = 534.5999755859375
$rngData[2,6] = 1061.300048828125
$rngData[2,7] = 36556.84924316406
$rngData[2,8] = 0.01982637516904765
$rngData[2,9] = 183.9243623717895
$rngData[3,0] = "2024-08-30"
$rngData[3,2] = 1933.349975585938
$rngData[3,3] = 4511.7998046875
$rngData[3,4] = 6132.10009765625
$rngData[3,5] = 538.7000122070312
$rngData[3,6] = 1065.599975585938
$rngData[3,7] = 36614.849609375
$rngData[3,8] = 0.001586580009265521
$rngData[3,9] = 184.2161730883455
$rngData[4,0] = "2024-09-02"
$rngData[4,2] = 1943.699951171875
$rngData[4,3] = 4553.75
$rngData[4,4] = 6156.0498046875
$rngData[4,5] = 538.4000244140625
$rngData[4,6] = 1050.949951171875
$rngData[4,7] = 36697.599609375
$rngData[4,8] = 0.002260012013781763
$rngData[4,9] = 184.632503852658
$rngData[5,0] = "2024-09-03"
$rngData[5,2] = 1964.5
$rngData[5,3] = 4521.0498046875
$rngData[5,4] = 6153.5
$rngData[5,5] = 532.4500122070312
$rngData[5,6] = 1068.800048828125
$rngData[5,7] = 36751.5
$rngData[5,8] = 0.001468771559958659
$rngData[5,9] = 184.9036868233608
$rngData[6,0] = "2024-09-04"
$rngData[6,2] = 1941.25
$rngData[6,3] = 4512.35009765625
$rngData[6,4] = 6145.7001953125
$rngData[6,5] = 536.0499877929688
$rngData[6,6] = 1056.199951171875
$rngData[6,7] = 36590.25
$rngData[6,8] = -0.004387576017305416
$rngData[6,9] = 184.0924078415432
$rngData[7,0] = "2024-09-05"
$rngData[7,2] = 1922.449951171875
$rngData[7,3] = 4479.25
$rngData[7,4] = 6071.2001953125
$rngData[7,5] = 519.1500244140625
$rngData[7,6] = 1074.900024414062
$rngData[7,7] = 36246.25036621094
$rngData[7,8] = -0.009401401569791475
$rngData[7,9] = 182.361681189475
$rngData[8,0] = "2024-09-06"
$rngData[8,2] = 1933.150024414062
$rngData[8,3] = 4475.9501953125
$rngData[8,4] = 6149.2998046875
$rngData[8,5] = 524.8499755859375
$rngData[8,6] = 1112.650024414062
$rngData[8,7] = 36628.40014648438
$rngData[8,8] = 0.01054315346863246
$rngData[8,9] = 184.2843483810535
$rngData[9,0] = "2024-09-09"
$rngData[9,2] = 1901.849975585938
$rngData[9,3] = 4456.75
$rngData[9,4] = 6165.39990234375
$rngData[9,5] = 520.5999755859375
$rngData[9,6] = 1077.550048828125
$rngData[9,7] = 36223.09973144531
$rngData[9,8] = -0.01106519567925938
$rngData[9,9] = 182.2452060055923
$rngData[10,0] = "2024-09-10"
$rngData[10,2] = 1894.650024414062
$rngData[10,3] = 4449.5498046875
$rngData[10,4] = 6146.60009765625
$rngData[10,5] = 514.8499755859375
$rngData[10,6] = 1091.0
$rngData[10,7] = 36152.14953613281
$rngData[10,8] = -0.00195870027243715
$rngData[10,9] = 181.8882422709388
$rngData[11,0] = "2024-09-11"
$rngData[11,2] = 1912.300048828125
$rngData[11,3] = 4507.85009765625
$rngData[11,4] = 6343.35009765625
$rngData[11,5] = 525.75
$rngData[11,6] = 1077.849975585938
$rngData[11,7] = 36618.80041503906
$rngData[11,8] = 0.01290797047738057
$rngData[11,9] = 184.2360503323547
$rngData[12,0] = "2024-09-12"
$rngData[12,2] = 1910.150024414062
$rngData[12,3] = 4479.35009765625
$rngData[12,4] = 6299.2998046875
$rngData[12,5] = 514.3499755859375
$rngData[12,6] = 1083.75
$rngData[12,7] = 36399.69982910156
$rngData[12,8] = -0.005983281359689682
$rngData[12,9] = 183.1337142066183
$rngData[13,0] = "2024-09-13"
$rngData[13,2] = 1950.449951171875
$rngData[13,3] = 4517.7001953125
$rngData[13,4] = 6392.35009765625
$rngData[13,5] = 530.0499877929688
$rngData[13,6] = 1089.699951171875
$rngData[13,7] = 36989.09985351562
$rngData[13,8] = 0.01619244189323883
$rngData[13,9] = 186.0990962326019
$rngData[14,0] = "2024-09-16"
$rngData[14,2] = 1944.099975585938
$rngData[14,3] = 4522.60009765625
$rngData[14,4] = 6416.2001953125
$rngData[14,5] = 550.5999755859375
$rngData[14,6] = 1094.650024414062
$rngData[14,7] = 37262.35009765625
$rngData[14,8] = 0.007387318026736299
$rngData[14,9] = 187.4738694409604
$rngData[15,0] = "2024-09-17"
$rngData[15,2] = 1950.25
$rngData[15,3] = 4513.25
$rngData[15,4] = 6423.4501953125
$rngData[15,5] = 551.9000244140625
$rngData[15,6] = 1080.300048828125
$rngData[15,7] = 37225.50073242188
$rngData[15,8] = -0.000988916832615256
$rngData[15,9] = 187.2884733757947
$rngData[16,0] = "2024-09-18"
$rngData[16,2] = 1952.550048828125
$rngData[16,3] = 4505.64990234375
$rngData[16,4] = 6455.75
$rngData[16,5] = 551.9000244140625
$rngData[16,6] = 1065.800048828125
$rngData[16,7] = 37181.6005859375
$rngData[16,8] = -0.0011793030481962
$rngData[16,9] = 187.0676035082506
$rngData[17,0] = "2024-09-19"
$rngData[17,2] = 1892.150024414062
$rngData[17,3] = 4346.14990234375
$rngData[17,4] = 6366.2998046875
$rngData[17,5] = 538.1500244140625
$rngData[17,6] = 1060.75
$rngData[17,7] = 36280.90002441406
$rngData[17,8] = -0.02422436224717267
$rngData[17,9] = 182.5360101161563
$rngData[18,0] = "2024-09-20"
$rngData[18,2] = 1894.199951171875
$rngData[18,3] = 4296.14990234375
$rngData[18,4] = 6377.14990234375
$rngData[18,5] = 533.3499755859375
$rngData[18,6] = 1114.699951171875
$rngData[18,7] = 36414.14892578125
$rngData[18,8] = 0.003672701098305774
$rngData[18,9] = 183.2064103209902
$rngData[19,0] = "2024-09-23"
$rngData[19,2] = 1905.75
$rngData[19,3] = 4284.89990234375
$rngData[19,4] = 6373.10009765625
$rngData[19,5] = 539.0999755859375
$rngData[19,6] = 1106.699951171875
$rngData[19,7] = 36474.34936523438
$rngData[19,8] = 0.001653215610663443
$rngData[19,9] = 183.5092900185065
$rngData[20,0] = "2024-09-24"
$rngData[20,2] = 1896.449951171875
$rngData[20,3] = 4268.5
$rngData[20,4] = 6326.10009765625
$rngData[20,5] = 534.9000244140625
$rngData[20,6] = 1098.5
$rngData[20,7] = 36256.65014648438
$rngData[20,8] = -0.005968556603164541
$rngData[20,9] = 182.4140044338245
$rngData[21,0] = "2024-09-25"
$rngData[21,2] = 1898.599975585938
$rngData[21,3] = 4271.2998046875
$rngData[21,4] = 6344.10009765625
$rngData[21,5] = 539.5499877929688
$rngData[21,6] = 1088.599975585938
$rngData[21,7] = 36297.29931640625
$rngData[21,8] = 0.001121150733938296
$rngData[21,9] = 182.6185180287761

$ws.Range("A615:J636").Value = $rngData

